$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start from a clean sheet (old "peroos/INICIO/FIN" login-style demo rows) ---
$ws.Cells.Clear()

# --- Column widths: A wider for names, B sized for the "Numero Nomina" values ---
$ws.Columns("A").ColumnWidth = 42.26
$ws.Columns("B").ColumnWidth = 16.09

# --- Header row (bold, centered) ---
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Número Nomina"
$ws.Range("C1").Value = "Password"
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").HorizontalAlignment = -4108

# --- Row 2: inicio user ---
$ws.Range("A2").Value = "inicio"
$ws.Range("B2").Formula = '="68657"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("C2").Value = 123456

# --- Row 3: fin user ---
$ws.Range("A3").Value = "fin"
$ws.Range("B3").Formula = '="68657"'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("C3").Value = 123456

# --- Rows 4-11: blank placeholder rows ready for more pasted-in CSV rows ---
$ws.Range("A4:B11").VerticalAlignment = -4160
$ws.Range("A2:B3").VerticalAlignment = -4160

$ws.Range("A4").Select()
